# Updates the "Price" (column D) and "Volume(1h)" (column E) figures in the
# cryptos list, reflecting refreshed data from the GitHub Actions scraper run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold numeric-looking text (e.g. "1.00", "7.00",
# "0.0000240") that must stay as literal text rather than being coerced to
# numbers (which would drop meaningful trailing zeros / reformat as
# scientific notation). Forcing the NumberFormat to Text ("@") before
# assigning the value keeps it as the exact literal string, matching the
# original inline-string cell contents.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.892.84"
$ws.Range("E2").Value = "  -0.28%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.078.47"
$ws.Range("E3").Value = "  -1.36%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.07"
$ws.Range("E5").Value = "  +0.30%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.29"
$ws.Range("E6").Value = "  -3.58%  "

$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.077.93"
$ws.Range("E8").Value = "  -1.23%  "

$ws.Range("E9").Value = "  -1.60%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.41"
$ws.Range("E10").Value = "  -0.37%  "

$ws.Range("E11").Value = "  -2.71%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.471"
$ws.Range("E12").Value = "  -1.80%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000240"
$ws.Range("E13").Value = "  -3.04%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.15"
$ws.Range("E14").Value = "  -3.08%  "

$ws.Range("E15").Value = "  -2.33%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.592.43"
$ws.Range("E16").Value = "  -1.20%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "66.832.83"
$ws.Range("E17").Value = "  -0.33%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.00"
$ws.Range("E18").Value = "  -1.83%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.082.00"
$ws.Range("E19").Value = "  -1.21%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.39"
$ws.Range("E20").Value = "  +1.08%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "484.52"
$ws.Range("E21").Value = "  +1.37%  "

$ws.Range("E22").Value = "  +0.22%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.687"
$ws.Range("E23").Value = "  -3.66%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.68"
$ws.Range("E24").Value = "  -1.51%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.83"
$ws.Range("E25").Value = "  -3.82%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.23"
$ws.Range("E26").Value = "  -3.36%  "

$ws.Range("E27").Value = "  +2.39%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.71"
$ws.Range("E29").Value = "  -3.44%  "

$ws.Range("E30").Value = "  -5.14%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.62"
$ws.Range("E31").Value = "  -1.65%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "27.75"
$ws.Range("E32").Value = "  -3.35%  "

$ws.Range("E33").Value = "  -2.10%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0₃0904"
$ws.Range("E34").Value = "  -7.17%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  -0.02%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.69"
$ws.Range("E36").Value = "  -2.94%  "

$ws.Range("E37").Value = "  -2.59%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "46.41"
$ws.Range("E38").Value = "  -2.78%  "

$ws.Range("E39").Value = "  +0.74%  "

$ws.Range("E40").Value = "  -5.30%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.301"
$ws.Range("E41").Value = "  -3.19%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.30"
$ws.Range("E42").Value = "  -3.70%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.769.13"
$ws.Range("E43").Value = "  -1.58%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "372.22"
$ws.Range("E44").Value = "  -2.32%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.52"
$ws.Range("E45").Value = "  -2.39%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0345"
$ws.Range("E46").Value = "  -3.13%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "135.55"
$ws.Range("E47").Value = "  -0.29%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "24.38"
$ws.Range("E49").Value = "  -1.86%  "

$ws.Range("E50").Value = "  -1.61%  "

$ws.Range("E51").Value = "  -3.12%  "

Write-Output "Updated cryptos list price/volume cells"
